$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 used to hold the "Thyroid" placeholder row; replace it with the
# new "Try" test-folder row. Cells are written in this order (C, B, A) so
# that the new shared-string entries land in the same order Excel recorded
# them: "Try" (C18), "Try folder" (B18), "Complete try" (A18).
$ws.Range("C18").Value = "Try"
$ws.Range("B18").Value = "Try folder"
$ws.Range("A18").Value = "Complete try"

$ws.Range("A18").Select()
